$p = $ppt.ActivePresentation

# ----------------------------------------------------------------------
# 1) Slide 2 ("Gliederung" / table of contents): add a new sub-bullet
#    "Netzwerk" (level 2 / lvl="1") right after the existing "Labyrinth"
#    sub-bullet in the content placeholder.
# ----------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$toc = $slide2.Shapes.Item(2).TextFrame.TextRange

# Append a new paragraph after the last one ("Labyrinth"); it inherits the
# same indent level / bullet formatting as the paragraph it follows.
$toc.InsertAfter("`r ")
$toc2 = $slide2.Shapes.Item(2).TextFrame.TextRange
$toc2.InsertAfter("Netzwerk")

# ----------------------------------------------------------------------
# 2) Slide 5 title: merge the two runs "Prototyp: " and "NETZWERK" into a
#    single run containing "Prototyp: NETZWERK".
# ----------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$title5 = $slide5.Shapes.Item(1).TextFrame.TextRange
$title5.Text = "placeholder"
$title5b = $slide5.Shapes.Item(1).TextFrame.TextRange
$title5b.Text = "Prototyp: NETZWERK"

# ----------------------------------------------------------------------
# 3) Slide 5 content: add a new sub-bullet "Player Collision (-> Kampfsystem)"
#    right after the existing "Disconnect-Handling" sub-bullet.
# ----------------------------------------------------------------------
$body5 = $slide5.Shapes.Item(2).TextFrame.TextRange
$body5.InsertAfter("`r ")
$body5b = $slide5.Shapes.Item(2).TextFrame.TextRange
$body5b.InsertAfter("Player ")
$body5c = $slide5.Shapes.Item(2).TextFrame.TextRange
$body5c.InsertAfter("Collision")
$body5d = $slide5.Shapes.Item(2).TextFrame.TextRange
$body5d.InsertAfter(" ")
$body5e = $slide5.Shapes.Item(2).TextFrame.TextRange
$body5e.InsertAfter("(-> Kampfsystem)")
